$d = $word.ActiveDocument

# The document carries three logo pictures inline in its headers/footers:
#   - two copies of the Pearson Edexcel "PearsonLogo.png" picture
#     (one in the "first page" footer, one in the "default" footer),
#     both currently named "image2.png" and which must be renamed to
#     "image1.png";
#   - one "BTec_Logo-Orange" picture in the "first page" header,
#     currently named "image1.jpg" and which must be renamed to
#     "image2.jpg".
#
# InlineShapes don't expose a settable .Name directly (same as real
# Word), so each picture is briefly promoted to a floating Shape -
# which does expose .Name - renamed, then converted back to an inline
# shape in place.

function Rename-LogoPicture($inlineShape) {
    $descr = $inlineShape.AlternativeText

    $newName = $null
    if ($descr -like "*PearsonLogo.png") {
        $newName = "image1.png"
    } elseif ($descr -eq "BTec_Logo-Orange") {
        $newName = "image2.jpg"
    }

    if ($newName -ne $null) {
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape() | Out-Null
    }
}

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $count = $hdr.Range.InlineShapes.Count
            for ($j = 1; $j -le $count; $j++) {
                Rename-LogoPicture $hdr.Range.InlineShapes.Item($j)
            }
        }

        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $count = $ftr.Range.InlineShapes.Count
            for ($j = 1; $j -le $count; $j++) {
                Rename-LogoPicture $ftr.Range.InlineShapes.Item($j)
            }
        }
    }
}
